# Change the table style applied to the table on slide 5
# (Slide "B1- TYPES OF FINANCIAL DOCUMENTS" -> the 3-column table)
# from the custom "Table_0" style to the built-in PowerPoint table
# style {E69BC347-AEEF-44CE-9E9A-F788FC60E955}.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# The table lives in the 2nd shape on the slide (a graphicFrame hosting
# an a:tbl). Find it defensively by checking HasTable in case shape
# ordering ever differs.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
        break
    }
}

$tbl = $tableShape.Table
$tbl.ApplyStyle("{E69BC347-AEEF-44CE-9E9A-F788FC60E955}")
